# Apply the update to the COASTAL_STATE_SOCIO_ECONOMIC sheet:
# Column F (rows 2-26) previously held a simple formula-driven ramp
# (10, 15, 20, ... +5 each row). It is replaced with actual observed
# values for most rows, and left blank for rows 2, 5, 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COASTAL_STATE_SOCIO_ECONOMIC")

# Map of row number -> new F value ($null means "clear the cell")
$values = @{
    2  = $null
    3  = 0.94
    4  = 0.77
    5  = $null
    6  = $null
    7  = $null
    8  = 0.47
    9  = 0.72
    10 = 0.85
    11 = 0.94
    12 = 0.69
    13 = 1.04
    14 = 1.1000000000000001
    15 = 1.1599999999999999
    16 = 0.64
    17 = 0.43
    18 = 0.77
    19 = 1.19
    20 = 2.38
    21 = 2.29
    22 = 0.7
    23 = 0.64
    24 = 0.78
    25 = 1.77
    26 = 0.76
}

foreach ($row in 2..26) {
    $cell = $ws.Cells.Item($row, 6)  # column F
    $val = $values[$row]
    if ($null -eq $val) {
        $cell.ClearContents()
    } else {
        $cell.Value = $val
    }
    # Column F switches from the red 2-decimal-place number style to the
    # plain "General" style already used by columns B/C (style index 10).
    $ws.Cells.Item($row, 2).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = $false

# The previously selected/active sheet switches to COASTAL_STATE_SOCIO_ECONOMIC
$ws.Activate()
$ws.Range("F15").Select()

$wb.Save()
